# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (between "2021-Q4" and "总计") that
#    mirrors the "2021-Q4" sheet's layout/format, with its own figures.
# 2) Prepend a "2022-Q1" row to the "总计" (totals) sheet, shifting the
#    existing rows down and renumbering the helper index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$srcQ4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $srcQ4)
$q1.Name = "2022-Q1"

# Copy header (with its formatting) and the two data rows' layout +
# formatting/shared values straight from 2021-Q4 - the sheets are
# structurally identical, only a handful of figures differ.
$srcQ4.Range("B1:H1").Copy($q1.Range("B1"))
$srcQ4.Range("A2:H3").Copy($q1.Range("A2"))

# The fund-scale / position figures are stored as text, not numbers, so
# enter them prefixed with an apostrophe in a scratch cell and copy only
# the *value* across - a plain .Value assignment would be auto-coerced
# to a number by Excel.
$scratch = $q1.Range("Z1")

function Set-TextValue($cell, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}

Set-TextValue $q1.Range("D2") "1.34"
Set-TextValue $q1.Range("E2") "82.00"
Set-TextValue $q1.Range("F2") "5.70"
Set-TextValue $q1.Range("G2") "0.0764"

Set-TextValue $q1.Range("D3") "1.34"
Set-TextValue $q1.Range("E3") "82.00"
Set-TextValue $q1.Range("F3") "5.70"
Set-TextValue $q1.Range("G3") "0.0764"

$scratch.Clear()

# ---------------------------------------------------------------------
# 2) "总计" sheet - insert the 2022-Q1 row at the top
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

# The inserted row inherits stray formatting from the header - strip it
# back to the plain (unstyled) look used by every other data row.
$total.Range("B2:D2").ClearFormats()

# Restore the index column's style (it's blank after the insert).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122) | Out-Null

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.15

# Renumber the helper index column (0-based) for every data row.
for ($i = 0; $i -le 5; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
}

# Leave the originally-active sheet selected (adding/naming sheets above
# shifts Excel's focus) so unrelated workbook view state doesn't change.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q1 sheet added; total sheet updated"
